# Applies "Various comments and sequence name changes" edit:
#  - Sheet "Tests" (sheet1): adds rows 5-20 of workflow/status pairs
#  - Sheet "Result" (sheet2): adds rows 2-20 of workflow/status/PASS data
#  - Makes "Result" the active/selected sheet, updates selections

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tests")
$ws2 = $wb.Worksheets.Item("Result")

# Data for rows 2..20, column A = workflow file, column B = status
$rows = @(
    @("Framework\InitAllSettings.xaml","Success"),
    @("Framework\InitAllApplications.xaml","Success"),
    @("Framework\CloseAllApplications.xaml","Success"),
    @("Framework\CloseAllApplications.xaml","AppEx"),
    @("Framework\InitAllSettings.xaml","Success"),
    @("Framework\InitAllSettings.xaml","Success"),
    @("Framework\InitAllApplications.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","Success"),
    @("Framework\CloseAllApplications.xaml","Success"),
    @("Test_Framework\Test_ProcessTransaction.xaml","AppEx"),
    @("Test_Framework\Test_ProcessTransaction.xaml","AppEx")
)

# --- Sheet1 "Tests": only rows 2-4 already have data; fill rows 5-20 (A,B) ---
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws1.Cells.Item($r, 2).Value = $rows[$i][1]
}

# --- Sheet2 "Result": fill rows 2-20 (A,B,C); C is always "PASS" ---
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws2.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws2.Cells.Item($r, 3).Value = "PASS"
}

# --- View/selection changes: Result becomes the active/selected tab ---
[void]$ws1.Range("A10").Select()
$ws2.Activate()
[void]$ws2.Range("A17").Select()
